# Refresh NATMI edge-weight metrics (columns G-T) for the Cthrc1-Ror2
# ligand-receptor sheet to match the new TPM computation. Columns A-F
# (cluster/gene labels, cell counts) are unaffected by the re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Cthrc1/Ror2)
$ws.Cells.Item(2, 7).Value = 0.8141179999999999  # column G
$ws.Cells.Item(2, 8).Value = 1.628236  # column H
$ws.Cells.Item(2, 9).Value = 0.03711595173929667  # column I
$ws.Cells.Item(2, 10).Value = 0.02558673725508277  # column J
$ws.Cells.Item(2, 13).Value = 0.071358  # column M
$ws.Cells.Item(2, 14).Value = 0.142716  # column N
$ws.Cells.Item(2, 15).Value = 0.01919591193090569  # column O
$ws.Cells.Item(2, 16).Value = 0.01411929935366186  # column P
$ws.Cells.Item(2, 17).Value = 0.058093832244  # column Q
$ws.Cells.Item(2, 18).Value = 0.232375328976  # column R
$ws.Cells.Item(2, 19).Value = 0.0007124745408192849  # column S
$ws.Cells.Item(2, 20).Value = 0.0003612668027880058  # column T

# Row 3: ECs -> FAPs (Cthrc1/Ror2)
$ws.Cells.Item(3, 7).Value = 0.8141179999999999  # column G
$ws.Cells.Item(3, 8).Value = 1.628236  # column H
$ws.Cells.Item(3, 9).Value = 0.03711595173929667  # column I
$ws.Cells.Item(3, 10).Value = 0.02558673725508277  # column J
$ws.Cells.Item(3, 15).Value = 0.7144019644080171  # column O
$ws.Cells.Item(3, 16).Value = 0.7882033865305114  # column P
$ws.Cells.Item(3, 17).Value = 2.16204096083  # column Q
$ws.Cells.Item(3, 18).Value = 12.97224576498  # column R
$ws.Cells.Item(3, 19).Value = 0.0265157088334267  # column S
$ws.Cells.Item(3, 20).Value = 0.02016755295472264  # column T

# Row 4: ECs -> MuSCs (Cthrc1/Ror2)
$ws.Cells.Item(4, 7).Value = 0.8141179999999999  # column G
$ws.Cells.Item(4, 8).Value = 1.628236  # column H
$ws.Cells.Item(4, 9).Value = 0.03711595173929667  # column I
$ws.Cells.Item(4, 10).Value = 0.02558673725508277  # column J
$ws.Cells.Item(4, 13).Value = 0.9728370000000001  # column M
$ws.Cells.Item(4, 14).Value = 1.945674  # column N
$ws.Cells.Item(4, 15).Value = 0.2617014683024538  # column O
$ws.Cells.Item(4, 16).Value = 0.1924910567184946  # column P
$ws.Cells.Item(4, 17).Value = 0.7920041127659999  # column Q
$ws.Cells.Item(4, 18).Value = 3.168016451064  # column R
$ws.Cells.Item(4, 19).Value = 0.009713299067616955  # column S
$ws.Cells.Item(4, 20).Value = 0.004925218092209356  # column T

# Row 5: ECs -> Neutrophils (Cthrc1/Ror2)
$ws.Cells.Item(5, 7).Value = 0.8141179999999999  # column G
$ws.Cells.Item(5, 8).Value = 1.628236  # column H
$ws.Cells.Item(5, 9).Value = 0.03711595173929667  # column I
$ws.Cells.Item(5, 10).Value = 0.02558673725508277  # column J
$ws.Cells.Item(5, 11).Value = 1  # column K
$ws.Cells.Item(5, 12).Value = 0.3333333333333333  # column L
$ws.Cells.Item(5, 13).Value = 0.017474  # column M
$ws.Cells.Item(5, 14).Value = 0.052422  # column N
$ws.Cells.Item(5, 15).Value = 0.004700655358623364  # column O
$ws.Cells.Item(5, 16).Value = 0.005186257397332197  # column P
$ws.Cells.Item(5, 17).Value = 0.014225897932  # column Q
$ws.Cells.Item(5, 18).Value = 0.085355387592  # column R
$ws.Cells.Item(5, 19).Value = 0.0001744692974337311  # column S
$ws.Cells.Item(5, 20).Value = 0.0001326994053627683  # column T

# Row 6: FAPs -> ECs (Cthrc1/Ror2)
$ws.Cells.Item(6, 9).Value = 0.9011867647895317  # column I
$ws.Cells.Item(6, 10).Value = 0.9318808175952528  # column J
$ws.Cells.Item(6, 13).Value = 0.071358  # column M
$ws.Cells.Item(6, 14).Value = 0.142716  # column N
$ws.Cells.Item(6, 15).Value = 0.01919591193090569  # column O
$ws.Cells.Item(6, 16).Value = 0.01411929935366186  # column P
$ws.Cells.Item(6, 17).Value = 1.410536178674  # column Q
$ws.Cells.Item(6, 18).Value = 8.463217072044  # column R
$ws.Cells.Item(6, 19).Value = 0.01729910177019767  # column S
$ws.Cells.Item(6, 20).Value = 0.01315750422556253  # column T

# Row 7: FAPs -> FAPs (Cthrc1/Ror2)
$ws.Cells.Item(7, 9).Value = 0.9011867647895317  # column I
$ws.Cells.Item(7, 10).Value = 0.9318808175952528  # column J
$ws.Cells.Item(7, 15).Value = 0.7144019644080171  # column O
$ws.Cells.Item(7, 16).Value = 0.7882033865305114  # column P
$ws.Cells.Item(7, 19).Value = 0.6438095950641471  # column S
$ws.Cells.Item(7, 20).Value = 0.7345116162713999  # column T

# Row 8: FAPs -> MuSCs (Cthrc1/Ror2)
$ws.Cells.Item(8, 9).Value = 0.9011867647895317  # column I
$ws.Cells.Item(8, 10).Value = 0.9318808175952528  # column J
$ws.Cells.Item(8, 13).Value = 0.9728370000000001  # column M
$ws.Cells.Item(8, 14).Value = 1.945674  # column N
$ws.Cells.Item(8, 15).Value = 0.2617014683024538  # column O
$ws.Cells.Item(8, 16).Value = 0.1924910567184946  # column P
$ws.Cells.Item(8, 17).Value = 19.230104325411  # column Q
$ws.Cells.Item(8, 18).Value = 115.380625952466  # column R
$ws.Cells.Item(8, 19).Value = 0.2358418995601586  # column S
$ws.Cells.Item(8, 20).Value = 0.1793787233146049  # column T

# Row 9: FAPs -> Neutrophils (Cthrc1/Ror2)
$ws.Cells.Item(9, 9).Value = 0.9011867647895317  # column I
$ws.Cells.Item(9, 10).Value = 0.9318808175952528  # column J
$ws.Cells.Item(9, 11).Value = 1  # column K
$ws.Cells.Item(9, 12).Value = 0.3333333333333333  # column L
$ws.Cells.Item(9, 13).Value = 0.017474  # column M
$ws.Cells.Item(9, 14).Value = 0.052422  # column N
$ws.Cells.Item(9, 15).Value = 0.004700655358623364  # column O
$ws.Cells.Item(9, 16).Value = 0.005186257397332197  # column P
$ws.Cells.Item(9, 17).Value = 0.3454091928886667  # column Q
$ws.Cells.Item(9, 18).Value = 3.108682735998  # column R
$ws.Cells.Item(9, 19).Value = 0.004236168395028366  # column S
$ws.Cells.Item(9, 20).Value = 0.004832973783685355  # column T

# Row 10: MuSCs -> ECs (Cthrc1/Ror2)
$ws.Cells.Item(10, 7).Value = 1.353296  # column G
$ws.Cells.Item(10, 8).Value = 2.706592  # column H
$ws.Cells.Item(10, 9).Value = 0.06169728347117155  # column I
$ws.Cells.Item(10, 10).Value = 0.04253244514966441  # column J
$ws.Cells.Item(10, 13).Value = 0.071358  # column M
$ws.Cells.Item(10, 14).Value = 0.142716  # column N
$ws.Cells.Item(10, 15).Value = 0.01919591193090569  # column O
$ws.Cells.Item(10, 16).Value = 0.01411929935366186  # column P
$ws.Cells.Item(10, 17).Value = 0.096568495968  # column Q
$ws.Cells.Item(10, 18).Value = 0.386273983872  # column R
$ws.Cells.Item(10, 19).Value = 0.001184335619888732  # column S
$ws.Cells.Item(10, 20).Value = 0.0006005283253113151  # column T

# Row 11: MuSCs -> FAPs (Cthrc1/Ror2)
$ws.Cells.Item(11, 7).Value = 1.353296  # column G
$ws.Cells.Item(11, 8).Value = 2.706592  # column H
$ws.Cells.Item(11, 9).Value = 0.06169728347117155  # column I
$ws.Cells.Item(11, 10).Value = 0.04253244514966441  # column J
$ws.Cells.Item(11, 15).Value = 0.7144019644080171  # column O
$ws.Cells.Item(11, 16).Value = 0.7882033865305114  # column P
$ws.Cells.Item(11, 17).Value = 3.59392788776  # column Q
$ws.Cells.Item(11, 18).Value = 21.56356732656  # column R
$ws.Cells.Item(11, 19).Value = 0.04407666051044324  # column S
$ws.Cells.Item(11, 20).Value = 0.03352421730438871  # column T

# Row 12: MuSCs -> MuSCs (Cthrc1/Ror2)
$ws.Cells.Item(12, 7).Value = 1.353296  # column G
$ws.Cells.Item(12, 8).Value = 2.706592  # column H
$ws.Cells.Item(12, 9).Value = 0.06169728347117155  # column I
$ws.Cells.Item(12, 10).Value = 0.04253244514966441  # column J
$ws.Cells.Item(12, 13).Value = 0.9728370000000001  # column M
$ws.Cells.Item(12, 14).Value = 1.945674  # column N
$ws.Cells.Item(12, 15).Value = 0.2617014683024538  # column O
$ws.Cells.Item(12, 16).Value = 0.1924910567184946  # column P
$ws.Cells.Item(12, 17).Value = 1.316536420752  # column Q
$ws.Cells.Item(12, 18).Value = 5.266145683008  # column R
$ws.Cells.Item(12, 19).Value = 0.01614626967467831  # column S
$ws.Cells.Item(12, 20).Value = 0.008187115311680314  # column T

# Row 13: MuSCs -> Neutrophils (Cthrc1/Ror2)
$ws.Cells.Item(13, 7).Value = 1.353296  # column G
$ws.Cells.Item(13, 8).Value = 2.706592  # column H
$ws.Cells.Item(13, 9).Value = 0.06169728347117155  # column I
$ws.Cells.Item(13, 10).Value = 0.04253244514966441  # column J
$ws.Cells.Item(13, 11).Value = 1  # column K
$ws.Cells.Item(13, 12).Value = 0.3333333333333333  # column L
$ws.Cells.Item(13, 13).Value = 0.017474  # column M
$ws.Cells.Item(13, 14).Value = 0.052422  # column N
$ws.Cells.Item(13, 15).Value = 0.004700655358623364  # column O
$ws.Cells.Item(13, 16).Value = 0.005186257397332197  # column P
$ws.Cells.Item(13, 17).Value = 0.023647494304  # column Q
$ws.Cells.Item(13, 18).Value = 0.141884965824  # column R
$ws.Cells.Item(13, 19).Value = 0.0002900176661612673  # column S
$ws.Cells.Item(13, 20).Value = 0.000220584208284073  # column T

# The refreshed run no longer includes a "Neutrophils" sending-cluster
# block (previously rows 14-17) - delete those rows and shift the rest up.
$ws.Range("A14:T17").EntireRow.Delete()
